# Amélioration du propagateur d'égalité.
# Replace the numeric ID values in column A (rows 2-9) of the first sheet
# with textual "LCP-n" identifiers, and move the active selection to B11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

for ($i = 1; $i -le 8; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = "LCP-$i"
}

$ws.Activate()
$ws.Range("B11").Select()
